# Update cryptos list: refresh prices and 1h volume percentages,
# and swap the Fetch.AI / USDe rows (47 & 48) with their new data.
# Numeric-looking price strings are forced to Text format before
# assignment (then reset to the default "Normal" style) so Excel
# stores them as plain text instead of auto-converting to numbers,
# matching the original inlineStr/text cell representation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.825.69"
$ws.Range("E2").Value = "  +4.27%  "
$ws.Range("D3").Value = "3.064.87"
$ws.Range("E3").Value = "  +2.37%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.47%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "3.049.15"
$ws.Range("E8").Value = "  +1.94%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.138"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +10.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.06%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.82"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.56%  "
$ws.Range("E15").Value = "  -0.01%  "
$ws.Range("D16").Value = "3.570.51"
$ws.Range("E16").Value = "  +2.37%  "
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "3.059.04"
$ws.Range("E18").Value = "  +2.15%  "
$ws.Range("D19").Value = "61.784.99"
$ws.Range("E19").Value = "  +4.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "448.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.729"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +3.77%  "
$ws.Range("E28").Value = "  +0.27%  "
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.08%  "
$ws.Range("D34").Value = "0.0₃0801"
$ws.Range("E34").Value = "  +2.85%  "
$ws.Range("E35").Value = "  +1.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.17"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.05"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.13%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "424.60"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.37%  "
$ws.Range("E42").Value = "  +4.71%  "
$ws.Range("D43").Value = "2.764.29"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.266"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "36.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.65%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.08"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.31%  "
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.36%  "
